$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2635626.5
$ws.Range("I86").Value = 3463
$ws.Range("J86").Value = 4789214.5
$ws.Range("K86").Value = 3463
$ws.Range("L86").Value = 4789214.5
$ws.Range("M86").Value = -2340
$ws.Range("N86").Value = -4791460.5
$ws.Range("H89").Value = 2635626.5
$ws.Range("I89").Value = 3463
$ws.Range("J89").Value = 4789214.5
$ws.Range("K89").Value = 17315
$ws.Range("L89").Value = 23946072.5
$ws.Range("M89").Value = -11699
$ws.Range("N89").Value = -23957304.5
$ws.Range("H116").Value = 9105.862999999999
$ws.Range("I116").Value = 4281.5835
$ws.Range("J116").Value = 14895
$ws.Range("K116").Value = 4281.5835
$ws.Range("L116").Value = 14895
$ws.Range("M116").Value = -839.5834999999997
$ws.Range("N116").Value = -21779
$ws.Range("H137").Value = 381433.34
$ws.Range("I137").Value = 247304
$ws.Range("J137").Value = 839708.5600000001
$ws.Range("K137").Value = 741912
$ws.Range("L137").Value = 2519125.68
$ws.Range("M137").Value = -739362
$ws.Range("N137").Value = -2524225.68
$ws.Range("H138").Value = 4511.5176
$ws.Range("I138").Value = 2243.121
$ws.Range("K138").Value = 6729.363
$ws.Range("M138").Value = -1589.363
$ws.Range("H141").Value = 2839.2666
$ws.Range("J141").Value = 7949.9
$ws.Range("L141").Value = 23849.7
$ws.Range("N141").Value = -34209.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 120.2
$ws.Range("I4").Value = 100.25
$ws.Range("K4").Value = 100.25
$ws.Range("M4").Value = 15.75
$ws.Range("H61").Value = 1796.5834
$ws.Range("I61").Value = 1521.3414
$ws.Range("K61").Value = 1521.3414
$ws.Range("M61").Value = -1309.3414
$ws.Range("H63").Value = 7609.1113
$ws.Range("I63").Value = 3994
$ws.Range("K63").Value = 3994
$ws.Range("M63").Value = -3308
$ws.Range("H66").Value = 7609.1113
$ws.Range("I66").Value = 3994
$ws.Range("K66").Value = 19970
$ws.Range("M66").Value = -16538
$ws.Range("H74").Value = 1506.2
$ws.Range("I74").Value = 1195.919
$ws.Range("J74").Value = 5333
$ws.Range("K74").Value = 1195.919
$ws.Range("L74").Value = 5333
$ws.Range("M74").Value = -321.9190000000001
$ws.Range("N74").Value = -7081
$ws.Range("H77").Value = 1506.2
$ws.Range("I77").Value = 1195.919
$ws.Range("J77").Value = 5333
$ws.Range("K77").Value = 5979.595
$ws.Range("L77").Value = 26665
$ws.Range("M77").Value = -1611.595
$ws.Range("N77").Value = -35401
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378
$ws.Range("H132").Value = 3732.4856
$ws.Range("I132").Value = 2506.6365
$ws.Range("J132").Value = 5807
$ws.Range("K132").Value = 7519.9095
$ws.Range("L132").Value = 17421
$ws.Range("M132").Value = -4989.9095
$ws.Range("N132").Value = -22481
$ws.Range("H136").Value = 1796.5834
$ws.Range("I136").Value = 1521.3414
$ws.Range("K136").Value = 4564.0242
$ws.Range("M136").Value = -2014.0242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4999
$ws.Range("I54").Value = 4999
$ws.Range("K54").Value = 4999
$ws.Range("M54").Value = -4515
$ws.Range("H86").Value = 1135137.1
$ws.Range("I86").Value = 1418700.1
$ws.Range("J86").Value = 885
$ws.Range("K86").Value = 1418700.1
$ws.Range("L86").Value = 885
$ws.Range("M86").Value = -1417577.1
$ws.Range("N86").Value = -3131
$ws.Range("H89").Value = 1135137.1
$ws.Range("I89").Value = 1418700.1
$ws.Range("J89").Value = 885
$ws.Range("K89").Value = 7093500.5
$ws.Range("L89").Value = 4425
$ws.Range("M89").Value = -7087884.5
$ws.Range("N89").Value = -15657
$ws.Range("H134").Value = 34815.848
$ws.Range("I134").Value = 4100.793
$ws.Range("K134").Value = 12302.379
$ws.Range("M134").Value = -9767.378999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 399.05713
$ws.Range("I7").Value = 387.30768
$ws.Range("K7").Value = 387.30768
$ws.Range("M7").Value = -274.30768
$ws.Range("H58").Value = 124228.57
$ws.Range("I58").Value = 155512.31
$ws.Range("J58").Value = 4614.294
$ws.Range("K58").Value = 155512.31
$ws.Range("L58").Value = 4614.294
$ws.Range("M58").Value = -155309.31
$ws.Range("N58").Value = -5020.294
$ws.Range("H105").Value = 4116.5
$ws.Range("J105").Value = 2933.6667
$ws.Range("L105").Value = 2933.6667
$ws.Range("N105").Value = -6427.6667
$ws.Range("H107").Value = 468.7143
$ws.Range("I107").Value = 421.83334
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 421.83334
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1498.16666
$ws.Range("N107").Value = -4590
$ws.Range("H122").Value = 4179.75
$ws.Range("I122").Value = 2849
$ws.Range("J122").Value = 5510.5
$ws.Range("K122").Value = 8547
$ws.Range("L122").Value = 16531.5
$ws.Range("M122").Value = -6097
$ws.Range("N122").Value = -21431.5
$ws.Range("H132").Value = 2225.4639
$ws.Range("I132").Value = 1887.4615
$ws.Range("K132").Value = 5662.3845
$ws.Range("M132").Value = -3132.3845
$ws.Range("H134").Value = 330084.53
$ws.Range("I134").Value = 190927.48
$ws.Range("K134").Value = 572782.4400000001
$ws.Range("M134").Value = -570247.4400000001
$ws.Range("H136").Value = 124228.57
$ws.Range("I136").Value = 155512.31
$ws.Range("J136").Value = 4614.294
$ws.Range("K136").Value = 466536.93
$ws.Range("L136").Value = 13842.882
$ws.Range("M136").Value = -463986.93
$ws.Range("N136").Value = -18942.882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1573764.4
$ws.Range("I132").Value = 501249.5
$ws.Range("J132").Value = 2002770.4
$ws.Range("K132").Value = 4511245.5
$ws.Range("L132").Value = 18024933.6
$ws.Range("M132").Value = -4508715.5
$ws.Range("N132").Value = -18029993.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 142.35
$ws.Range("I2").Value = 117
$ws.Range("J2").Value = 218.4
$ws.Range("K2").Value = 117
$ws.Range("L2").Value = 218.4
$ws.Range("M2").Value = -4
$ws.Range("N2").Value = -444.4
$ws.Range("H97").Value = 902.92
$ws.Range("I97").Value = 691.6875
$ws.Range("K97").Value = 691.6875
$ws.Range("M97").Value = -195.6875
$ws.Range("H126").Value = 2971.7222
$ws.Range("I126").Value = 1957.6666
$ws.Range("K126").Value = 5872.9998
$ws.Range("M126").Value = -3402.9998
$ws.Range("H132").Value = 252305.73
$ws.Range("I132").Value = 265603.28
$ws.Range("J132").Value = 168087.83
$ws.Range("K132").Value = 796809.8400000001
$ws.Range("L132").Value = 504263.49
$ws.Range("M132").Value = -794279.8400000001
$ws.Range("N132").Value = -509323.49

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 44065.207
$ws.Range("I40").Value = 52309.8
$ws.Range("J40").Value = 2842.25
$ws.Range("K40").Value = 52309.8
$ws.Range("L40").Value = 2842.25
$ws.Range("M40").Value = -52173.8
$ws.Range("N40").Value = -3114.25
$ws.Range("H132").Value = 1684.2609
$ws.Range("I132").Value = 889.5472
$ws.Range("K132").Value = 2668.6416
$ws.Range("M132").Value = -138.6415999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 492.9091
$ws.Range("I113").Value = 398.72726
$ws.Range("K113").Value = 1196.18178
$ws.Range("M113").Value = 973.8182200000001
$ws.Range("H132").Value = 18491.732
$ws.Range("I132").Value = 1000.65955
$ws.Range("J132").Value = 81728.69500000001
$ws.Range("K132").Value = 3001.97865
$ws.Range("L132").Value = 245186.085
$ws.Range("M132").Value = -471.97865
$ws.Range("N132").Value = -250246.085
